# Fruta / hortaliza, semanal
# Insert the latest week's two rows (1a nueva(o) / 2a nueva(o), origen Peru)
# at the top of this product's data block (row 245), pushing all existing
# rows for this block down by 2 (245-357 -> 247-359).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 245, shifting 245:357 down to 247:359.
$ws.Rows("245:246").Insert()

# New row 245: "1a nueva(o)"
$ws.Range("A245").Value = 8
$ws.Range("B245").Value = "Terminal La Palmera de La Serena"
$ws.Range("C245").Value = "Coquimbo"
$ws.Range("D245").Value = (Get-Date -Year 2021 -Month 9 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E245").Value = 4
$ws.Range("F245").Value = 100112045
$ws.Range("G245").Value = "Zapallo"
$ws.Range("H245").Value = "Camote"
$ws.Range("I245").Value = "1a nueva(o)"
$ws.Range("J245").Value = 800
$ws.Range("K245").Value = 950
$ws.Range("L245").Value = 1000
$ws.Range("M245").Value = 975
$ws.Range("N245").Value = '$/kilo (volumen en unidades)'
$ws.Range("O245").Value = "Perú"
$ws.Range("P245").Value = 975
$ws.Range("Q245").Value = 1
$ws.Range("R245").Value = "Hortaliza"

# New row 246: "2a nueva(o)"
$ws.Range("A246").Value = 8
$ws.Range("B246").Value = "Terminal La Palmera de La Serena"
$ws.Range("C246").Value = "Coquimbo"
$ws.Range("D246").Value = (Get-Date -Year 2021 -Month 9 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E246").Value = 4
$ws.Range("F246").Value = 100112045
$ws.Range("G246").Value = "Zapallo"
$ws.Range("H246").Value = "Camote"
$ws.Range("I246").Value = "2a nueva(o)"
$ws.Range("J246").Value = 560
$ws.Range("K246").Value = 850
$ws.Range("L246").Value = 900
$ws.Range("M246").Value = 875
$ws.Range("N246").Value = '$/kilo (volumen en unidades)'
$ws.Range("O246").Value = "Perú"
$ws.Range("P246").Value = 875
$ws.Range("Q246").Value = 1
$ws.Range("R246").Value = "Hortaliza"
